# Update Fgf7-Nrp1 LR-pairs sheet with refreshed TPM numbers and a new
# "Resolving-Mac" sending/target cluster (rows 2-9 updated in place,
# rows 10-13 appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data for data rows 2..13 (columns A..T), in final/target order.
$rows = @(
    @("FAPs","Fgf7","Nrp1","ECs",3,1,11.14069866666667,33.422096,0.9461861633339901,0.94618616333399,3,1,133.7780026666667,401.334008,0.50863533211804,0.5086353321180399,1490.380415937863,13413.42374344077,0.4812637134328781,0.4812637134328779),
    @("FAPs","Fgf7","Nrp1","FAPs",3,1,11.14069866666667,33.422096,0.9461861633339901,0.94618616333399,3,1,52.441971,157.325913,0.1993888292903622,0.1993888292903622,584.2401963970721,5258.161767573649,0.1886589513979037,0.1886589513979037),
    @("FAPs","Fgf7","Nrp1","MuSCs",3,1,11.14069866666667,33.422096,0.9461861633339901,0.94618616333399,3,1,21.197691,63.593073,0.08059542216956049,0.08059542216956046,236.157087860112,2125.413790741008,0.07625827328489965,0.07625827328489962),
    @("FAPs","Fgf7","Nrp1","Resolving-Mac",3,1,11.14069866666667,33.422096,0.9461861633339901,0.94618616333399,3,1,55.59592133333333,166.787764,0.2113804164220374,0.2113804164220373,619.3774066703716,5574.396660033344,0.2000052252183087,0.2000052252183086),
    @("MuSCs","Fgf7","Nrp1","ECs",2,0.6666666666666666,0.5338349999999999,1.601505,0.04533892403128163,0.04533892403128163,3,1,133.7780026666667,401.334008,0.50863533211804,0.5086353321180399,71.41538005356,642.7384204820401,0.02306097868252551,0.02306097868252551),
    @("MuSCs","Fgf7","Nrp1","FAPs",2,0.6666666666666666,0.5338349999999999,1.601505,0.04533892403128163,0.04533892403128163,3,1,52.441971,157.325913,0.1993888292903622,0.1993888292903622,27.995359588785,251.958236299065,0.009040074983881914,0.009040074983881912),
    @("MuSCs","Fgf7","Nrp1","MuSCs",2,0.6666666666666666,0.5338349999999999,1.601505,0.04533892403128163,0.04533892403128163,3,1,21.197691,63.593073,0.08059542216956049,0.08059542216956046,11.316069374985,101.844624374865,0.003654109723014774,0.003654109723014773),
    @("MuSCs","Fgf7","Nrp1","Resolving-Mac",2,0.6666666666666666,0.5338349999999999,1.601505,0.04533892403128163,0.04533892403128163,3,1,55.59592133333333,166.787764,0.2113804164220374,0.2113804164220373,29.67904866497999,267.11143798482,0.009583760641859427,0.009583760641859425),
    @("Resolving-Mac","Fgf7","Nrp1","ECs",1,0.3333333333333333,0.09978633333333332,0.299359,0.008474912634728231,0.008474912634728231,3,1,133.7780026666667,401.334008,0.50863533211804,0.5086353321180399,13.34921636676355,120.142947300872,0.004310640002636368,0.004310640002636366),
    @("Resolving-Mac","Fgf7","Nrp1","FAPs",1,0.3333333333333333,0.09978633333333332,0.299359,0.008474912634728231,0.008474912634728231,3,1,52.441971,157.325913,0.1993888292903622,0.1993888292903622,5.232991998863,47.096927989767,0.001689802908576561,0.001689802908576561),
    @("Resolving-Mac","Fgf7","Nrp1","MuSCs",1,0.3333333333333333,0.09978633333333332,0.299359,0.008474912634728231,0.008474912634728231,3,1,21.197691,63.593073,0.08059542216956049,0.08059542216956046,2.115239860023,19.037158740207,0.000683039161646064,0.0006830391616460638),
    @("Resolving-Mac","Fgf7","Nrp1","Resolving-Mac",1,0.3333333333333333,0.09978633333333332,0.299359,0.008474912634728231,0.008474912634728231,3,1,55.59592133333333,166.787764,0.2113804164220374,0.2113804164220373,5.547713138141777,49.92941824327599,0.001791430561869239,0.001791430561869239)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    for ($c = 1; $c -le $data.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
